$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (column E / F headers) ---
$ws.Range("E1").Value = "memories"
$ws.Range("F1").Value = "quantity"

# --- Data value updates (serial numbers replacing memory module refs) ---
$ws.Range("E2").Value = "SN001"
$ws.Range("E3").Value = "SN002"
$ws.Range("E5").Value = "SN004,SN005"
$ws.Range("E7").Value = "SN005,SN006"
$ws.Range("E11").Value = "SN003"

# --- Font change: memories column values now rendered in Times New Roman ---
$ws.Range("E2").Font.Name = "Times New Roman"
$ws.Range("E3").Font.Name = "Times New Roman"
$ws.Range("E5").Font.Name = "Times New Roman"
$ws.Range("E7").Font.Name = "Times New Roman"
$ws.Range("E11").Font.Name = "Times New Roman"
$ws.Range("E7").WrapText = $true

# --- E4 had no value; its leftover formatting is cleared back to workbook default ---
$ws.Range("E4").Font.Name = "Arial"
$ws.Range("E4").Font.Size = 10
$ws.Range("E4").VerticalAlignment = -4107

# --- Quantity column (F) values are now center aligned ---
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F3").HorizontalAlignment = -4108
$ws.Range("F4").HorizontalAlignment = -4108
$ws.Range("F5").HorizontalAlignment = -4108
$ws.Range("F6").HorizontalAlignment = -4108
$ws.Range("F7").HorizontalAlignment = -4108
$ws.Range("F8").HorizontalAlignment = -4108
$ws.Range("F9").HorizontalAlignment = -4108
$ws.Range("F10").HorizontalAlignment = -4108
$ws.Range("F11").HorizontalAlignment = -4108

# --- prototype_reference column (G) + D11 style refresh ---
$ws.Range("G2:G11").Locked = $true
$ws.Range("D11").Locked = $true

# --- Row heights collapse back to the sheet default now that wrapped text fits ---
$ws.Rows(3).RowHeight = 12.8
$ws.Rows(5).RowHeight = 12.8
$ws.Rows(7).RowHeight = 12.8

# --- Column G default style / active selection ---
$ws.Range("F2").Select()
